$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '''26.997.14'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '''  -2.46%  '
$ws.Range('E2').Style = 'Normal'

# Row 3
$ws.Range('D3').Value = '''1.816.37'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '''  -1.53%  '
$ws.Range('E3').Style = 'Normal'

# Row 4
$ws.Range('D4').Value = '''1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '''  -1.06%  '
$ws.Range('E4').Style = 'Normal'

# Row 5
$ws.Range('D5').Value = '''310.69'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '''  -2.63%  '
$ws.Range('E5').Style = 'Normal'

# Row 6
$ws.Range('D6').Value = '''1.000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '''  -1.04%  '
$ws.Range('E6').Style = 'Normal'

# Row 7
$ws.Range('D7').Value = '''0.4217'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '''  -2.09%  '
$ws.Range('E7').Style = 'Normal'

# Row 8
$ws.Range('D8').Value = '''0.3666'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '''  -1.94%  '
$ws.Range('E8').Style = 'Normal'

# Row 9
$ws.Range('D9').Value = '''0.07206'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '''  -1.78%  '
$ws.Range('E9').Style = 'Normal'

# Row 10
$ws.Range('D10').Value = '''0.8393'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '''  -4.18%  '
$ws.Range('E10').Style = 'Normal'

# Row 11
$ws.Range('D11').Value = '''20.74'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '''  -3.89%  '
$ws.Range('E11').Style = 'Normal'

# Row 12
$ws.Range('D12').Value = '''1.809.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '''  -1.65%  '
$ws.Range('E12').Style = 'Normal'

# Row 13
$ws.Range('D13').Value = '''6.628'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '''  -1.38%  '
$ws.Range('E13').Style = 'Normal'

# Row 14
$ws.Range('D14').Value = '''0.07063'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '''  -0.75%  '
$ws.Range('E14').Style = 'Normal'

# Row 15
$ws.Range('D15').Value = '''5.277'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '''  -2.98%  '
$ws.Range('E15').Style = 'Normal'

# Row 16
$ws.Range('D16').Value = '''88.90'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '''  +0.22%  '
$ws.Range('E16').Style = 'Normal'

# Row 17
$ws.Range('D17').Value = '''1.003'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '''  -1.14%  '
$ws.Range('E17').Style = 'Normal'

# Row 18
$ws.Range('D18').Value = '''0.000008773'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '''  -2.22%  '
$ws.Range('E18').Style = 'Normal'

# Row 19
$ws.Range('E19').Value = '''  -1.03%  '
$ws.Range('E19').Style = 'Normal'

# Row 20
$ws.Range('B20').Value = '''Avalanche'
$ws.Range('B20').Style = 'Normal'
$ws.Range('C20').Value = '''https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('C20').Style = 'Normal'
$ws.Range('D20').Value = '''14.92'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '''  -3.50%  '
$ws.Range('E20').Style = 'Normal'

# Row 21
$ws.Range('B21').Value = '''WrappedBTC'
$ws.Range('B21').Style = 'Normal'
$ws.Range('C21').Value = '''https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('C21').Style = 'Normal'
$ws.Range('D21').Value = '''27.046.95'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '''  -2.28%  '
$ws.Range('E21').Style = 'Normal'

# Row 22
$ws.Range('B22').Value = '''Uniswap'
$ws.Range('B22').Style = 'Normal'
$ws.Range('C22').Value = '''https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('C22').Style = 'Normal'
$ws.Range('D22').Value = '''5.112'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '''  -1.87%  '
$ws.Range('E22').Style = 'Normal'

# Row 23
$ws.Range('B23').Value = '''Cosmos'
$ws.Range('B23').Style = 'Normal'
$ws.Range('C23').Value = '''https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('C23').Style = 'Normal'
$ws.Range('D23').Value = '''10.80'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '''  -2.60%  '
$ws.Range('E23').Style = 'Normal'

# Row 24
$ws.Range('B24').Value = '''WrappedliquidstakedEther2.0'
$ws.Range('B24').Style = 'Normal'
$ws.Range('C24').Value = '''https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('C24').Style = 'Normal'
$ws.Range('D24').Value = '''2.039.70'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '''  -1.69%  '
$ws.Range('E24').Style = 'Normal'

# Row 25
$ws.Range('B25').Value = '''Toncoin'
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = '''https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = '''1.974'
$ws.Range('D25').Style = 'Normal'

# Row 26
$ws.Range('B26').Value = '''Monero'
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = '''https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = '''151.57'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '''  -2.56%  '
$ws.Range('E26').Style = 'Normal'

# Row 27
$ws.Range('B27').Value = '''LidoDAOToken'
$ws.Range('B27').Style = 'Normal'
$ws.Range('C27').Value = '''https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('C27').Style = 'Normal'
$ws.Range('D27').Value = '''2.227'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '''  +2.87%  '
$ws.Range('E27').Style = 'Normal'

# Row 28
$ws.Range('B28').Value = '''EthereumClassic'
$ws.Range('B28').Style = 'Normal'
$ws.Range('C28').Value = '''https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('C28').Style = 'Normal'
$ws.Range('D28').Value = '''18.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '''  -1.97%  '
$ws.Range('E28').Style = 'Normal'

# Row 29
$ws.Range('B29').Value = '''InternetComputer(DFINITY)'
$ws.Range('B29').Style = 'Normal'
$ws.Range('C29').Value = '''https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('C29').Style = 'Normal'
$ws.Range('D29').Value = '''5.204'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '''  -2.99%  '
$ws.Range('E29').Style = 'Normal'

# Row 30
$ws.Range('B30').Value = '''BitcoinCash'
$ws.Range('B30').Style = 'Normal'
$ws.Range('C30').Value = '''https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('C30').Style = 'Normal'
$ws.Range('D30').Value = '''115.87'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '''  -2.42%  '
$ws.Range('E30').Style = 'Normal'

# Row 31
$ws.Range('B31').Value = '''Stellar'
$ws.Range('B31').Style = 'Normal'
$ws.Range('C31').Value = '''https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('C31').Style = 'Normal'
$ws.Range('D31').Value = '''0.08751'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '''  -2.05%  '
$ws.Range('E31').Style = 'Normal'

# Row 32
$ws.Range('B32').Value = '''ARBITRUM'
$ws.Range('B32').Style = 'Normal'
$ws.Range('C32').Value = '''https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('C32').Style = 'Normal'
$ws.Range('D32').Value = '''1.171'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '''  -4.78%  '
$ws.Range('E32').Style = 'Normal'

# Row 33
$ws.Range('B33').Value = '''HuobiToken'
$ws.Range('B33').Style = 'Normal'
$ws.Range('C33').Value = '''https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('C33').Style = 'Normal'
$ws.Range('D33').Value = '''2.960'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '''  +2.62%  '
$ws.Range('E33').Style = 'Normal'

# Row 34
$ws.Range('B34').Value = '''ImmutableX'
$ws.Range('B34').Style = 'Normal'
$ws.Range('C34').Value = '''https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('C34').Style = 'Normal'
$ws.Range('D34').Value = '''0.7368'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '''  -4.90%  '
$ws.Range('E34').Style = 'Normal'

# Row 35
$ws.Range('B35').Value = '''Filecoin'
$ws.Range('B35').Style = 'Normal'
$ws.Range('C35').Value = '''https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('C35').Style = 'Normal'
$ws.Range('D35').Value = '''4.401'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '''  -3.15%  '
$ws.Range('E35').Style = 'Normal'

# Row 36
$ws.Range('B36').Value = '''Frax'
$ws.Range('B36').Style = 'Normal'
$ws.Range('C36').Value = '''https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('C36').Style = 'Normal'
$ws.Range('D36').Value = '''1.000'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '''  -1.18%  '
$ws.Range('E36').Style = 'Normal'

# Row 37
$ws.Range('B37').Value = '''TrustWalletToken'
$ws.Range('B37').Style = 'Normal'
$ws.Range('C37').Value = '''https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('C37').Style = 'Normal'
$ws.Range('D37').Value = '''1.087'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '''  -4.13%  '
$ws.Range('E37').Style = 'Normal'

# Row 38
$ws.Range('B38').Value = '''VeChain'
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = '''https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = '''0.01950'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '''  -1.09%  '
$ws.Range('E38').Style = 'Normal'

# Row 39
$ws.Range('B39').Value = '''Hedera'
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = '''https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = '''0.05221'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '''  -2.08%  '
$ws.Range('E39').Style = 'Normal'

# Row 40
$ws.Range('B40').Value = '''FraxShare'
$ws.Range('B40').Style = 'Normal'
$ws.Range('C40').Value = '''https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('C40').Style = 'Normal'
$ws.Range('D40').Value = '''7.266'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '''  -0.30%  '
$ws.Range('E40').Style = 'Normal'

# Row 41
$ws.Range('B41').Value = '''MXToken'
$ws.Range('B41').Style = 'Normal'
$ws.Range('C41').Value = '''https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('C41').Style = 'Normal'
$ws.Range('D41').Value = '''2.863'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '''  -2.06%  '
$ws.Range('E41').Style = 'Normal'

# Row 42
$ws.Range('B42').Value = '''Algorand'
$ws.Range('B42').Style = 'Normal'
$ws.Range('C42').Value = '''https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('C42').Style = 'Normal'
$ws.Range('D42').Value = '''0.1684'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '''  +0.03%  '
$ws.Range('E42').Style = 'Normal'

# Row 43
$ws.Range('B43').Value = '''TheSandbox'
$ws.Range('B43').Style = 'Normal'
$ws.Range('C43').Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('C43').Style = 'Normal'
$ws.Range('D43').Value = '''0.5021'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '''  -1.79%  '
$ws.Range('E43').Style = 'Normal'

# Row 44
$ws.Range('B44').Value = '''Aptos'
$ws.Range('B44').Style = 'Normal'
$ws.Range('C44').Value = '''https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('C44').Style = 'Normal'
$ws.Range('D44').Value = '''8.593'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '''  -2.26%  '
$ws.Range('E44').Style = 'Normal'

# Row 45
$ws.Range('B45').Value = '''EnergySwap'
$ws.Range('B45').Style = 'Normal'
$ws.Range('C45').Value = '''https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('C45').Style = 'Normal'
$ws.Range('D45').Value = '''10.49'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '''  -1.49%  '
$ws.Range('E45').Style = 'Normal'

# Row 46
$ws.Range('B46').Value = '''Quant'
$ws.Range('B46').Style = 'Normal'
$ws.Range('C46').Value = '''https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('C46').Style = 'Normal'
$ws.Range('D46').Value = '''106.06'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '''  -2.75%  '
$ws.Range('E46').Style = 'Normal'

# Row 47
$ws.Range('B47').Value = '''Decentraland'
$ws.Range('B47').Style = 'Normal'
$ws.Range('C47').Value = '''https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('C47').Style = 'Normal'
$ws.Range('D47').Value = '''0.4704'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '''  -0.75%  '
$ws.Range('E47').Style = 'Normal'

# Row 48
$ws.Range('B48').Value = '''PaxDollar'
$ws.Range('B48').Style = 'Normal'
$ws.Range('C48').Value = '''https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('C48').Style = 'Normal'
$ws.Range('D48').Value = '''1.000'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '''  -1.19%  '
$ws.Range('E48').Style = 'Normal'

# Row 49
$ws.Range('B49').Value = '''Cronos'
$ws.Range('B49').Style = 'Normal'
$ws.Range('C49').Value = '''https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('C49').Style = 'Normal'
$ws.Range('D49').Value = '''0.06339'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '''  -2.08%  '
$ws.Range('E49').Style = 'Normal'

# Row 50
$ws.Range('B50').Value = '''NEARProtocol'
$ws.Range('B50').Style = 'Normal'
$ws.Range('C50').Value = '''https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('C50').Style = 'Normal'
$ws.Range('D50').Value = '''1.639'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '''  -3.00%  '
$ws.Range('E50').Style = 'Normal'

# Row 51
$ws.Range('D51').Value = '''1.869'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '''  +1.35%  '
$ws.Range('E51').Style = 'Normal'
